$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Resolving-Mac" -> "ECs" and updated numeric values
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il19"
$ws.Range("C2").Value = "Il20ra"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.5
$ws.Range("G2").Value = 0.0621735
$ws.Range("H2").Value = 0.124347
$ws.Range("I2").Value = 0.400116482236459
$ws.Range("J2").Value = 0.3077956989247312
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02041766666666667
$ws.Range("N2").Value = 0.061253
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.0012694377985
$ws.Range("R2").Value = 0.007616626791
$ws.Range("S2").Value = 0.400116482236459
$ws.Range("T2").Value = 0.3077956989247312

# Row 3: new row for "Inflammatory-Mac"
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("B3").Value = "Il19"
$ws.Range("C3").Value = "Il20ra"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.09321499999999999
$ws.Range("H3").Value = 0.279645
$ws.Range("I3").Value = 0.5998835177635411
$ws.Range("J3").Value = 0.6922043010752689
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.02041766666666667
$ws.Range("N3").Value = 0.061253
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.001903232798333333
$ws.Range("R3").Value = 0.017129095185
$ws.Range("S3").Value = 0.5998835177635411
$ws.Range("T3").Value = 0.6922043010752689
